$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E3").Value = 17
$ws.Range("E4").Value = 16
$ws.Range("E16").Value = 262
$ws.Range("E18").Value = 72
